# AWS DEVICE FARM.pptx edit script
# 1. Delete the blank slide at position 8 (sldId 261) from the deck.
# 2. Merge the three split runs in slide 4's "As the tests complete..."
#    paragraph back into a single run.
# 3. Turn the (now-last) trailing blank slide into a centred "Thank you"
#    slide: drop its empty content placeholder, resize/recentre the title
#    and set its text.

$p = $ppt.ActivePresentation

# --- 1. Remove the blank slide that used to sit at index 8 ------------
$p.Slides.Item(8).Delete()

# --- 2. Fix up the "As the tests complete..." paragraph on slide 4 ----
$slide4 = $p.Slides.Item(4)
$contentShape = $slide4.Shapes.Item(2)
$para = $contentShape.TextFrame.TextRange.Paragraphs(3)
# Work around the no-op optimisation when the assigned text already
# equals the paragraph's current (concatenated) text: stage a throwaway
# value first so the run-merge actually takes effect.
$para.Text = "__tmp__"
$para.Text = "As the tests complete, test report with results, logs, screenshots and performance will be updated."

# --- 3. Turn the final blank slide into the "Thank you" slide ---------
$lastSlide = $p.Slides.Item($p.Slides.Count)

$title = $lastSlide.Shapes.Item(1)
$title.Name = "Title 3"
$title.Left = 51.333388
$title.Top = 208.0
$title.Width = 676.902985
$title.Height = 104.0
$title.TextFrame.TextRange.Text = "Thank you"
$title.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$content = $lastSlide.Shapes.Item(2)
$content.Delete()
